$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) column F changed on both the "展览"
# sheet and the mirrored "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 73
    $ws.Range("F3").Value = 1440
    $ws.Range("F4").Value = 13
    $ws.Range("F5").Value = 21
    $ws.Range("F8").Value = 224
}
